# Revert "adding term 2.0.0"
# - drop the extra "Include from FSIII 12" sheet that was added
# - restore the previous Metadata values (Version/Date/Contact)

$wb = $excel.ActiveWorkbook

# Remove the last "Include from FSIII 12" worksheet entirely.
[void]$wb.Worksheets("Include from FSIII 12").Delete()

# Restore the prior metadata property values on the Metadata sheet.
$ws = $wb.Worksheets("Metadata")
$ws.Range("B3").Value = "1.1.0"
$ws.Range("B8").Value = "2023-07-10T23:08:03+02:00"
$ws.Range("B10").Value = "No display for ContactDetail"
